# "Counts instead of binary features"
#
# 1. B69's note gets extended.
# 2. A parallel "WITH SUBSET AND WORD COUNTS" block of results is added in
#    column H, mirroring the existing "WITH THE SUBSET OF FEATURES" block
#    that lives in column B (rows 69-100).
# 3. The active selection moves to H74 (where the new work was happening).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the note in B69 -------------------------------------------
$ws.Range("B69").Value = "WITH THE SUBSET OF FEATURES (OCCURS 40:100 times), note avg length didn't add anything"

# --- 2. Mirror the B69:B100 formatting onto H69:H100 ----------------------
# (copies the Courier-New / left-aligned style used throughout the block,
# including the blank separator rows, without inventing new style records)
$ws.Range("B69:B100").Copy()
$ws.Range("H69:H100").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# H69 uses the Courier-New style (like the rest of the block) even though
# B69 itself (the note row) is unstyled, so fix its style up explicitly.
$ws.Range("B70").Copy()
$ws.Range("H69").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Populate the new "WITH SUBSET AND WORD COUNTS" block in column H --
$ws.Range("H69").Value = "WITH SUBSET AND WORD COUNTS"
$ws.Range("H70").Value = "Folds: 10"

$ws.Range("H72").Value = "           Gaussian       MNB Bernoulli"
$ws.Range("H73").Value = "Accuracy   0.789474  0.815789  0.684211"
$ws.Range("H74").Value = "Precision  0.783784  0.823529  0.793103"
$ws.Range("H75").Value = "Recall            1  0.965517  0.793103"
$ws.Range("H76").Value = "F1         0.878788  0.888889  0.793103"

$ws.Range("H78").Value = "           Gaussian       MNB Bernoulli"
$ws.Range("H79").Value = "Accuracy   0.368421  0.605263  0.631579"
$ws.Range("H80").Value = "Precision  0.461538   0.62963  0.681818"
$ws.Range("H81").Value = "Recall     0.545455  0.772727  0.681818"
$ws.Range("H82").Value = "F1              0.5  0.693878  0.681818"

$ws.Range("H84").Value = "           Gaussian       MNB Bernoulli"
$ws.Range("H85").Value = "Accuracy        0.5  0.578947  0.631579"
$ws.Range("H86").Value = "Precision       0.5  0.545455       0.6"
$ws.Range("H87").Value = "Recall     0.631579  0.947368  0.789474"
$ws.Range("H88").Value = "F1          0.55814  0.692308  0.681818"

$ws.Range("H90").Value = "           Gaussian       MNB Bernoulli"
$ws.Range("H91").Value = "Accuracy   0.763158  0.815789  0.710526"
$ws.Range("H92").Value = "Precision  0.763158  0.823529  0.821429"
$ws.Range("H93").Value = "Recall            1  0.965517  0.793103"
$ws.Range("H94").Value = "F1         0.865672  0.888889  0.807018"

$ws.Range("H96").Value = "           Gaussian       MNB Bernoulli"
$ws.Range("H97").Value = "Accuracy   0.631579  0.605263  0.421053"
$ws.Range("H98").Value = "Precision  0.631579  0.642857       0.4"
$ws.Range("H99").Value = "Recall     0.631579  0.473684  0.315789"
$ws.Range("H100").Value = "F1         0.631579  0.545455  0.352941"

# --- 4. Move the view / selection to where the edits happened -------------
$ws.Activate()
$win = $excel.ActiveWindow
try { $win.ScrollRow = 67 } catch {}
try { $win.ScrollColumn = 2 } catch {}
$ws.Range("H74").Select()

# --- 5. Best-effort: shrink the saved window width (cosmetic bookView) ----
try { $excel.ActiveWindow.Width = 13470 } catch {}
try { $excel.Width = 13470 } catch {}
